$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 270
    $ws.Range("F3").Value = 72
    $ws.Range("F5").Value = 6839
    $ws.Range("F6").Value = 5485
    $ws.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202409/RqtCRIaH1726800618582.jpeg"
    $ws.Range("F7").Value = 452
    $ws.Range("F11").Value = 241

    if ($name -eq "展览") {
        $ws.Range("F12").Value = 134
    } else {
        $ws.Range("F14").Value = 134
    }
}
